$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update RTM content cells (rows 6-29) to match the revised RTM ---
$ws.Range('B8').Value = 'User can donate online'
$ws.Range('F8').Value = 'Donation Module'
$ws.Range('G8').Value = 'FTC_Do12'
$ws.Range('B9').Value = 'User can shop online'
$ws.Range('F9').Value = 'Shopping Cart Module'
$ws.Range('G9').Value = 'FTC_SoCa8,'
$ws.Range('A10').Value = ''
$ws.Range('G10').Value = 'FTC_SoCa9'
$ws.Range('A11').Value = 'R4'
$ws.Range('B11').Value = 'Any user can view the foundation''s programs'
$ws.Range('C11').Value = 'Lowy Needed'
$ws.Range('D11').Value = 'Completed'
$ws.Range('E11').Value = 'Low'
$ws.Range('F11').Value = 'User Interface'
$ws.Range('G11').Value = 'UITC_Ev1'
$ws.Range('A12').Value = 'R5'
$ws.Range('B12').Value = 'User can cancel his/her donation'
$ws.Range('C12').Value = 'Highly Needed'
$ws.Range('D12').Value = 'In Progress'
$ws.Range('E12').Value = 'High'
$ws.Range('F12').Value = 'Donation Page'
$ws.Range('G12').Value = 'FTC_Do11'
$ws.Range('A13').Value = 'R6'
$ws.Range('B13').Value = 'User can cancel his/her order'
$ws.Range('C13').Value = 'Highly Needed'
$ws.Range('D13').Value = 'In Progress'
$ws.Range('E13').Value = 'Low'
$ws.Range('F13').Value = 'Shopping Cart Module'
$ws.Range('G13').Value = 'FTC_SoCa10'
$ws.Range('B14').Value = 'Admin can view records of volunteers'
$ws.Range('C14').Value = 'Highly Needed'
$ws.Range('D14').Value = 'In Progress'
$ws.Range('E14').Value = 'High'
$ws.Range('F14').Value = 'Back End Volunteer Module'
$ws.Range('G14').Value = 'UITC_ReAp1'
$ws.Range('B15').Value = 'Admin can view records of subscribers'
$ws.Range('F15').Value = 'Back End Subscriber Module'
$ws.Range('G15').Value = 'UITC_SeEm1'
$ws.Range('B16').Value = 'Admin can email volunteers'
$ws.Range('E16').Value = 'High'
$ws.Range('F16').Value = 'Back End Volunteer Module'
$ws.Range('G16').Value = 'UITC_ReAp2'
$ws.Range('B17').Value = 'Admin can email subscribers'
$ws.Range('G17').Value = 'UITC_SeEm2'
$ws.Range('B18').Value = 'User can select buy products in the shop'
$ws.Range('F18').Value = 'Shopping Cart Module'
$ws.Range('G18').Value = 'FTC_SoCa7,'
$ws.Range('A19').Value = ''
$ws.Range('B19').Value = ''
$ws.Range('C19').Value = ''
$ws.Range('D19').Value = ''
$ws.Range('E19').Value = ''
$ws.Range('F19').Value = ''
$ws.Range('G19').Value = 'FTC_SoCa6,'
$ws.Range('A20').Value = ''
$ws.Range('B20').Value = ''
$ws.Range('C20').Value = ''
$ws.Range('D20').Value = ''
$ws.Range('E20').Value = ''
$ws.Range('F20').Value = ''
$ws.Range('G20').Value = 'FTC_SoCa5'
$ws.Range('A21').Value = 'R12'
$ws.Range('B21').Value = 'Users can pay without PayPal account'
$ws.Range('G21').Value = 'FTC_SoCa8'
$ws.Range('A22').Value = 'R13'
$ws.Range('B22').Value = 'Users can view all ordered items'
$ws.Range('G22').Value = 'FTC_SoCa14'
$ws.Range('A23').Value = 'R14'
$ws.Range('B23').Value = 'Online transactions are secured'
$ws.Range('F23').Value = 'Shopping Cart and Donation Module'
$ws.Range('G23').Value = 'FTC_SoCa12,'
$ws.Range('A24').Value = ''
$ws.Range('B24').Value = ''
$ws.Range('C24').Value = ''
$ws.Range('D24').Value = ''
$ws.Range('E24').Value = ''
$ws.Range('F24').Value = ''
$ws.Range('G24').Value = 'FTC_SoCa13,'
$ws.Range('A25').Value = ''
$ws.Range('B25').Value = ''
$ws.Range('C25').Value = ''
$ws.Range('D25').Value = ''
$ws.Range('E25').Value = ''
$ws.Range('F25').Value = ''
$ws.Range('G25').Value = 'FTC_SoCa16-'
$ws.Range('G26').Value = 'FTC_SoCa23,'
$ws.Range('G27').Value = 'FTC_Do2-FTC_Do10'
$ws.Range('G28').Value = ''
$ws.Range('G29').Value = ''

# --- Clear the stray explicit-font formatting on G23/G24/G27 (now plain) ---
$ws.Range('G23').ClearFormats()
$ws.Range('G24').ClearFormats()
$ws.Range('G27').ClearFormats()

# --- Move the active selection (no more scrolled topLeftCell) ---
$ws.Range('G14').Select()
